$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# --- Rename header cells: "_old" -> "_FV2210", "_new" -> "_FV2304" ---
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2210"
        $newVal = $newVal -replace "_new$", "_FV2304"
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# --- Turn the data range into a proper Excel Table (ListObject) ---
# Stash the header row's current formatting in a scratch range, then clear the
# header formatting before creating the table. Excel's ListObjects.Add()
# otherwise "preserves" any pre-existing explicit header formatting as a
# table-level header dxf (headerRowDxfId) - which the target workbook does not
# have. After the table exists we restore the original header formatting.
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A1000:U1000")
$hdr.Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats

$hdr.ClearFormats()

$range = $ws.Range("A1:U66")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes, $null)
$table.Name = "Table1"
$table.TableStyle = ""

# restore original header formatting
$scratch.Copy()
$hdr.PasteSpecial(-4122) # xlPasteFormats
$scratch.ClearFormats()
$scratch.ClearContents()
$excel.CutCopyMode = $false

# --- Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
